# Add 2022-Q4 data:
#  - insert a new worksheet "2022-Q4" (positioned right after "总计", before
#    the existing "2022-Q3" sheet) by duplicating the "2022-Q3" sheet so it
#    inherits the same layout/styling, then overwrite its data with the new
#    quarter's figures.
#  - update the "总计" (totals) sheet with a new row for 2022-Q4 and shift
#    the existing 2022-Q3 / 2022-Q2 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet to get a same-styled "2022-Q4" sheet,
#    inserted right before it (i.e. right after "总计").
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Fill in the new "2022-Q4" sheet with the quarter's fund data
#    (3 funds this quarter, one more row than the template had).
# ---------------------------------------------------------------------

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "007497"
$q4.Range("C2").Value = "中庚价值灵动灵活配置混合"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "34.74"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "93.96"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "2.55"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.8859"
$q4.Range("H2").Value = 10

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "168105"
$q4.Range("C3").Value = "九泰泰富灵活配置混合（LOF）A"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.43"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "94.68"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "4.65"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0200"
$q4.Range("H3").Value = 2

# Row 4 is new - copy formatting from row 3 first, then set values.
$q4.Range("A3:H3").Copy()
$q4.Range("A4:H4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "015688"
$q4.Range("C4").Value = "九泰泰富灵活配置混合（LOF）C"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.00"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "94.68"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "4.65"
$q4.Range("G4").NumberFormat = "General"
$q4.Range("G4").Value = 0
$q4.Range("H4").Value = 2

# ---------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: new 2022-Q4 row, shift others down.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Add row 4 first (copy formatting from row 3, which already has the style
# used by the data rows), then overwrite row values top-down.
$totals.Range("A3:D3").Copy()
$totals.Range("A4:D4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 4 <- old row 3 (2022-Q2)
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q2"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0.01

# Row 3 <- old row 2 (2022-Q3)
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.02

# Row 2 <- new 2022-Q4 figures
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 0.91

# Restore the originally-selected tab ("2022-Q2") as the active sheet.
$wb.Worksheets.Item("2022-Q2").Activate()
